$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 560, shifting existing rows
# (560-665) down to (562-667), matching the target dataset which now has
# 667 data rows instead of 665.
$ws.Rows("560:561").Insert()

# --- New row 560 (new weekly record) ---
$ws.Range("A560").Value = 6
$ws.Range("B560").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C560").Value = "Metropolitana"
$ws.Range("D560").Value = "2023-02-27"
$ws.Range("E560").Value = 13
$ws.Range("F560").Value = 100112030
$ws.Range("G560").Value = "Poroto granado"
$ws.Range("H560").Value = "Sin especificar"
$ws.Range("I560").Value = "Primera"
$ws.Range("J560").Value = 1100
$ws.Range("K560").Value = 28000
$ws.Range("L560").Value = 30000
$ws.Range("M560").Value = 28909
$ws.Range("N560").Value = "`$/saco 25 kilos"
$ws.Range("O560").Value = "Región Metropolitana"
$ws.Range("P560").Value = 1156
$ws.Range("Q560").Value = 25
$ws.Range("R560").Value = "Hortaliza"

# --- New row 561 (new weekly record) ---
$ws.Range("A561").Value = 6
$ws.Range("B561").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C561").Value = "Metropolitana"
$ws.Range("D561").Value = "2023-02-27"
$ws.Range("E561").Value = 13
$ws.Range("F561").Value = 100112030
$ws.Range("G561").Value = "Poroto granado"
$ws.Range("H561").Value = "Sin especificar"
$ws.Range("I561").Value = "Primera"
$ws.Range("J561").Value = 370
$ws.Range("K561").Value = 28000
$ws.Range("L561").Value = 30000
$ws.Range("M561").Value = 28811
$ws.Range("N561").Value = "`$/saco 25 kilos"
$ws.Range("O561").Value = "Región del Maule"
$ws.Range("P561").Value = 1152
$ws.Range("Q561").Value = 25
$ws.Range("R561").Value = "Hortaliza"
